$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top; the existing rows 1-56 shift down
# to become rows 3-58 (old row 1, the text header row, lands on row 3;
# old row 2, the first data row, lands on row 4; and so on).
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# The freshly inserted row 1 starts out with no formatting. Copy the
# header formatting (bold, centered, bordered - style that travelled down
# with the old row 1 content onto row 3) onto the new row 1 so it keeps
# matching visual treatment.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 1: sequential numbers 0-13.
$headerValues = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13)
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerValues[$i]
}

# New row 2: blank except for E2 = "Washer".
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(2, $col).Value = ""
}
$ws.Cells.Item(2, 5).Value = "Washer"

# New row 3 (previously row 1) keeps its original header labels, but the
# thread_size / material_surface labels in columns M and N are cleared.
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(3, 14).Value = ""
